$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "2026-01-28 04:05"
$ws.Range("B4").Value = 39
$ws.Range("C4").Value = 7
